$d = $word.ActiveDocument

# The document contains four "<id>XXXX</id>" markers that were each split
# across three separate runs (one run for the literal "<id>" text, one for
# the inner id value, one for the closing "</id>" text). The edit merges
# each trio into a single run holding the full "<id>XXXX</id>" text, using
# the formatting of the first ("<id>") run - which Word's Find/Replace
# naturally does when the found range spans multiple runs: the replacement
# text inherits the formatting of the start of the found range and the
# runs it used to occupy collapse into one.

$ids = @("p132r_2", "p132v_1", "p132v_2", "p132v_3")

foreach ($id in $ids) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
